# MISS02TP001 template fix (#89)
# - Adds a "Type Day" legend: two new header columns on Sheet1
#   (Deploy Program / Deplot For User) with sample values, plus a
#   widened C:E header band.
# - Adds a second worksheet "Descrition" holding the Type-Day code
#   legend (English + Thai) that Sheet1's "Type Day" column refers to.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: new header cells -------------------------------------------
# Write content in the same order the original workbook's shared-string
# table grew in: D1, E1, D2, D5, then C1 (C1 reuses/replaces the old
# "Type day" slot).
$ws1.Range("D1").Value = "Deploy Program"
$ws1.Range("E1").Value = "Deplot For User"
$ws1.Range("D2").Value = "a"
$ws1.Range("D5").Value = "s"
$ws1.Range("C1").Value = "Type Day"

# Match formatting: D1:E1 use the same header style as C1, D2/D5 use the
# same body style as the rest of column A/C.
$ws1.Range("C1").Copy() | Out-Null
$ws1.Range("D1:E1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy() | Out-Null
$ws1.Range("D2").PasteSpecial(-4122) | Out-Null
$ws1.Range("D5").PasteSpecial(-4122) | Out-Null

$ws1.Application.CutCopyMode = $false

# Widen the new header columns C:E.
$ws1.Range("C1:E1").ColumnWidth = 12.67

$ws1.Range("C14").Select() | Out-Null

# --- New worksheet "Descrition" (Type-Day legend) -----------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Descrition"

$ws2.Range("B6").Value  = "H    "
$ws2.Range("D6").Value  = "Holiday"
$ws2.Range("E6").Value  = "วันหยุด"

$ws2.Range("B7").Value  = "W    "
$ws2.Range("D7").Value  = "Saturday Working "
$ws2.Range("E7").Value  = "วันทำงานที่ได้รับมอบหมาย"

$ws2.Range("B8").Value  = "S    "
$ws2.Range("D8").Value  = "Send Program Package"
$ws2.Range("E8").Value  = "วันส่งไฟล์แพคเกจ"

$ws2.Range("B9").Value  = "I    "
$ws2.Range("D9").Value  = "Deployment on PROD/QA."
$ws2.Range("E9").Value  = "วันส่งมอบบนเครื่อง PROD/QA."

$ws2.Range("B10").Value = "D    "
$ws2.Range("D10").Value = "V-Smart infrom issue no deployment"
$ws2.Range("E10").Value = "วันที่วีสมาร์ทส่งมอบ"

$ws2.Range("J8").Select() | Out-Null
